$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.52
$ws.Range("C2").Value = 16.940000000000001

$ws.Range("B3").Value = 18.7
$ws.Range("C3").Value = 7.91

$ws.Range("B4").Value = 17.329999999999998
$ws.Range("C4").Value = 64.349999999999994

$ws.Range("B5").Value = 17.16
$ws.Range("C5").Value = 85.87

$ws.Range("B6").Value = 17.55
$ws.Range("C6").Value = 93.34

$ws.Range("G6").Value = 2839.7139999999999

$ws.Range("G3").Select()
